$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 35 and 36 swapped rank (Fetch.AI moved above Maker) with updated price/volume
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '3.14'
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '3.732.58'
$ws.Range("E36").Value = '  +2.14%  '

# Price / Volume(1h) updates for remaining rows
$ws.Range("D2").Value = '69.581.16'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '3.503.24'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").Value = '195.80'
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -5.45%  '
$ws.Range("D10").Value = '0.643'
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").Value = '53.01'
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("E12").Value = '  -3.84%  '
$ws.Range("E13").Value = '  -1.37%  '
$ws.Range("D14").Value = '4.061.10'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '596.77'
$ws.Range("E15").Value = '  -2.65%  '
$ws.Range("D16").Value = '69.683.50'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '18.90'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = '3.498.40'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("E22").Value = '  +4.15%  '
$ws.Range("E23").Value = '  +3.39%  '
$ws.Range("D24").Value = '101.71'
$ws.Range("E24").Value = '  -4.78%  '
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("D26").Value = '3.09'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("D28").Value = '9.44'
$ws.Range("E28").Value = '  -2.74%  '
$ws.Range("D29").Value = '32.91'
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").Value = '4.28'
$ws.Range("E30").Value = '  +8.87%  '
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("D34").Value = '63.14'
$ws.Range("E34").Value = '  -0.48%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '0.0₃0805'
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("D39").Value = '3.62'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("D41").Value = '35.90'
$ws.Range("E41").Value = '  -2.25%  '
$ws.Range("D42").Value = '493.30'
$ws.Range("E42").Value = '  -4.08%  '
$ws.Range("D43").Value = '0.132'
$ws.Range("E43").Value = '  -3.83%  '
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("E45").Value = '  -2.86%  '
$ws.Range("E46").Value = '  -4.00%  '
$ws.Range("D47").Value = '3.27'
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("E49").Value = '  -4.20%  '
$ws.Range("D50").Value = '0.000242'
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Value = '128.22'
$ws.Range("E51").Value = '  -2.69%  '
